$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.334.68"
$ws.Range("E2").Value = "  -0.88%  "

$ws.Range("D3").Value = "3.752.60"
$ws.Range("E3").Value = "  -2.21%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.02"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.21%  "

$ws.Range("D7").Value = "3.750.24"
$ws.Range("E7").Value = "  -2.28%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.48"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("E13").Value = "  +6.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.06%  "

$ws.Range("D15").Value = "4.386.25"
$ws.Range("E15").Value = "  -2.02%  "

$ws.Range("D16").Value = "3.769.36"
$ws.Range("E16").Value = "  -1.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.65"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.89%  "

$ws.Range("D18").Value = "67.418.93"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.68%  "

$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.66"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.720"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.68%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.77"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -9.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.17"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("E28").Value = "  +3.28%  "

$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.90"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.93%  "

$ws.Range("D31").Value = "3.907.02"
$ws.Range("E31").Value = "  -1.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.69"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.45"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.93%  "

$ws.Range("E34").Value = "  -3.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.12"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.44%  "

$ws.Range("D36").Value = "3.718.95"
$ws.Range("E36").Value = "  -2.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.81"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.53%  "

$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.87"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.994"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.51%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.85"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "399.11"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000269"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -9.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.44"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0354"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.86%  "
